# Gantt.xlsx update:
#  - rename Sheet1 -> "Before covid"
#  - add a new worksheet "After Covid" (after "Before covid", becomes active)
#  - populate the new sheet with a second (unfinished) gantt-style chart
#  - move the selection around on sheet "Before covid" per the saved view

$wb = $excel.ActiveWorkbook

# --- sheet 1: rename, keep data as-is, just move the saved selection -------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Before covid"

# --- add the second sheet right after "Before covid" -----------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "After Covid"

# column widths (best effort, engine quantizes to 1/6 character steps)
$ws2.Columns.Item(1).ColumnWidth = 13
$ws2.Columns.Item(2).ColumnWidth = 15.5
$ws2.Columns.Item(3).ColumnWidth = 23.66666666666667
$ws2.Columns.Item(4).ColumnWidth = 17.33333333333333
$ws2.Columns.Item(5).ColumnWidth = 15.16666666666667
$ws2.Columns.Item(9).ColumnWidth = 11.66666666666667

# date header row (same style/number-format used for the dates on sheet 1)
$dateFmt = $ws1.Range("B1").NumberFormat
$ws2.Range("A2").Value = "date"
$ws2.Range("B2").Value = 43934
$ws2.Range("C2").Value = 43941
$ws2.Range("D2").Value = 43948
$ws2.Range("E2").Value = 43955
$ws2.Range("F2").Value = 43962
$ws2.Range("G2").Value = 43969
$ws2.Range("H2").Value = 43976
$ws2.Range("I2").Value = 43979
$ws2.Range("B2:I2").NumberFormat = $dateFmt

# row labels (column A) - typed in this order so new shared-strings line up
$ws2.Range("A3").Value = "Ideation"
$ws2.Range("A4").Value = "Solution"
$ws2.Range("A5").Value = "Baseline"
$ws2.Range("A11").Value = "Iterative"
$ws2.Range("A6").Value = "Testing"
$ws2.Range("A7").Value = "Management"
$ws2.Range("A8").Value = "Mathematical"
$ws2.Range("A9").Value = "API"
$ws2.Range("A10").Value = "GUI"

# "final touches" marker column (I) for every task row
$ws2.Range("I3").Value = "final touches"
$ws2.Range("I4").Value = "final touches"
$ws2.Range("I5").Value = "final touches"
$ws2.Range("I6").Value = "final touches"
$ws2.Range("I7").Value = "final touches"
$ws2.Range("I8").Value = "final touches"
$ws2.Range("I9").Value = "final touches"
$ws2.Range("I10").Value = "final touches"
$ws2.Range("I11").Value = "final touches"

# Baseline sub-tasks
$ws2.Range("B5").Value = "find new dataset"
$ws2.Range("C5").Value = "setup detection in one file"
$ws2.Range("D5").Value = "setup environment"
$ws2.Range("E5").Value = "setup training"

# Mathematical sub-tasks
$ws2.Range("C8").Value = "evaluate detection"
$ws2.Range("E8").Value = "evaluate training"

# page setup - portrait, no explicit paper size
$ws2.PageSetup.Orientation = 1

# --- selections / active views ---------------------------------------------
$null = $ws1.Range("F48").Select()
$null = $ws2.Activate()
$null = $ws2.Range("C10").Select()
